# Update crypto price/volume cells to refreshed values (GitHub Actions scrape update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.154.57"
$ws.Range("E2").Value = "  +3.35%  "
$ws.Range("D3").Value = "2.425.00"
$ws.Range("E3").Value = "  +4.42%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.97"
$ws.Range("E5").Value = "  +2.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.45"
$ws.Range("E6").Value = "  +6.25%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  +2.99%  "
$ws.Range("D9").Value = "2.425.31"
$ws.Range("E9").Value = "  +4.43%  "
$ws.Range("E10").Value = "  +5.80%  "
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.38"
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.353"
$ws.Range("E13").Value = "  +5.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.28"
$ws.Range("E14").Value = "  +8.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000174"
$ws.Range("E15").Value = "  +10.79%  "
$ws.Range("D16").Value = "2.820.77"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").Value = "62.000.67"
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("D18").Value = "2.427.27"
$ws.Range("E18").Value = "  +4.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.15"
$ws.Range("E19").Value = "  +6.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.19"
$ws.Range("E20").Value = "  +3.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.48"
$ws.Range("E21").Value = "  +3.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.78"
$ws.Range("E22").Value = "  +5.08%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.97"
$ws.Range("E24").Value = "  +4.00%  "
$ws.Range("E25").Value = "  +5.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.14"
$ws.Range("E26").Value = "  +11.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "559.04"
$ws.Range("E27").Value = "  +16.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.513.56"
$ws.Range("E29").Value = "  +3.19%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.35"
$ws.Range("E30").Value = "  +7.31%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0932"
$ws.Range("E31").Value = "  +10.13%  "
$ws.Range("E32").Value = "  +8.03%  "
$ws.Range("E33").Value = "  +4.32%  "
$ws.Range("E34").Value = "  +5.59%  "
$ws.Range("E35").Value = "  +4.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.76"
$ws.Range("E36").Value = "  +13.31%  "
$ws.Range("E37").Value = "  +15.82%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.82"
$ws.Range("E39").Value = "  +7.46%  "
$ws.Range("E40").Value = "  +4.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.76"
$ws.Range("E41").Value = "  +2.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "146.49"
$ws.Range("E42").Value = "  +3.73%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.28"
$ws.Range("E44").Value = "  +13.41%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "149.85"
$ws.Range("E45").Value = "  +6.96%  "
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0541"
$ws.Range("E47").Value = "  +7.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.38"
$ws.Range("E48").Value = "  +8.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.592"
$ws.Range("E49").Value = "  +5.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0227"
$ws.Range("E50").Value = "  +4.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0910"
$ws.Range("E51").Value = "  +2.60%  "
